$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA8").Value = -1
$ws.Range("AB8").Value = 0.8500000000000001
$ws.Range("AC8").Value = -1
$ws.Range("AD8").Value = 0.8999999999999999
$ws.Range("B8").Value = 7126858
$ws.Range("E8").Value = "Sava Strmec"
$ws.Range("F8").Value = "Lucko"
$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = "A"
$ws.Range("L8").Value = 2.75
$ws.Range("M8").Value = 3.4
$ws.Range("N8").Value = 2.2
$ws.Range("O8").Value = 3.75
$ws.Range("P8").Value = 3.6
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 0.5
$ws.Range("S8").Value = 1.95
$ws.Range("T8").Value = 1.85
$ws.Range("V8").Value = 1.9
$ws.Range("W8").Value = 1.9
$ws.Range("X8").Value = -1
$ws.Range("Z8").Value = 0.8
$ws.Range("AA9").Value = 0
$ws.Range("AB9").Value = 0
$ws.Range("AD9").Value = 0.8500000000000001
$ws.Range("B9").Value = 7126860
$ws.Range("E9").Value = "NK Tomislav"
$ws.Range("F9").Value = "NK Oriolik Oriovac"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = "D"
$ws.Range("L9").Value = 2.5
$ws.Range("N9").Value = 2.4
$ws.Range("O9").Value = 2.625
$ws.Range("P9").Value = 3.4
$ws.Range("Q9").Value = 2.3
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 1.975
$ws.Range("T9").Value = 1.725
$ws.Range("U9").Value = 3
$ws.Range("V9").Value = 1.95
$ws.Range("W9").Value = 1.85
$ws.Range("Y9").Value = 2.4
$ws.Range("Z9").Value = -1
$ws.Range("AA10").Value = 0.925
$ws.Range("AB10").Value = -1
$ws.Range("AC10").Value = 0.825
$ws.Range("AD10").Value = -1
$ws.Range("B10").Value = 7126857
$ws.Range("E10").Value = "NK Bistra"
$ws.Range("F10").Value = "NK Vrapce"
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 3
$ws.Range("K10").Value = "H"
$ws.Range("L10").Value = 2.875
$ws.Range("M10").Value = 3.5
$ws.Range("N10").Value = 2.1
$ws.Range("O10").Value = 2.9
$ws.Range("P10").Value = 3.5
$ws.Range("Q10").Value = 2.1
$ws.Range("R10").Value = 0.25
$ws.Range("S10").Value = 1.925
$ws.Range("T10").Value = 1.875
$ws.Range("U10").Value = 2.75
$ws.Range("V10").Value = 1.825
$ws.Range("W10").Value = 1.975
$ws.Range("X10").Value = 1.9
$ws.Range("Y10").Value = -1
$ws.Range("F12").Value = "Sava Strmec"
$ws.Range("F13").Value = "NK Tomislav"
$ws.Range("E19").Value = "Sava Strmec"
$ws.Range("F19").Value = "NK Vrapce"
$ws.Range("E21").Value = "NK Bistra"
$ws.Range("AA22").Value = 0.8
$ws.Range("AB22").Value = -1
$ws.Range("AC22").Value = 0.825
$ws.Range("AD22").Value = -1
$ws.Range("B22").Value = 7202435
$ws.Range("E22").Value = "NK Udarnik Kurilovec"
$ws.Range("F22").Value = "NK Mladost Petrinja"
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 2
$ws.Range("M22").Value = 3.4
$ws.Range("N22").Value = 3.1
$ws.Range("O22").Value = 2
$ws.Range("P22").Value = 3.4
$ws.Range("Q22").Value = 3.1
$ws.Range("R22").Value = -0.25
$ws.Range("S22").Value = 1.8
$ws.Range("T22").Value = 2
$ws.Range("X22").Value = 1
$ws.Range("AA23").Value = 0.35
$ws.Range("AB23").Value = -0.5
$ws.Range("AC23").Value = 0.4125
$ws.Range("AD23").Value = -0.5
$ws.Range("B23").Value = 7202436
$ws.Range("E23").Value = "Lucko"
$ws.Range("F23").Value = "NK Tondach"
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1
$ws.Range("L23").Value = 1.615
$ws.Range("M23").Value = 3.75
$ws.Range("N23").Value = 4.333
$ws.Range("O23").Value = 1.533
$ws.Range("P23").Value = 4
$ws.Range("Q23").Value = 5
$ws.Range("R23").Value = -0.75
$ws.Range("S23").Value = 1.7
$ws.Range("T23").Value = 2.1
$ws.Range("U23").Value = 2.75
$ws.Range("V23").Value = 1.825
$ws.Range("W23").Value = 1.975
$ws.Range("X23").Value = 0.5329999999999999
$ws.Range("AA24").Value = 0.8500000000000001
$ws.Range("AC24").Value = 0.8
$ws.Range("B24").Value = 7202437
$ws.Range("E24").Value = "NK Maksimir"
$ws.Range("F24").Value = "Sava Strmec"
$ws.Range("L24").Value = 1.4
$ws.Range("M24").Value = 4.333
$ws.Range("N24").Value = 6
$ws.Range("O24").Value = 1.4
$ws.Range("P24").Value = 4.333
$ws.Range("Q24").Value = 6
$ws.Range("R24").Value = -1.25
$ws.Range("S24").Value = 1.85
$ws.Range("T24").Value = 1.95
$ws.Range("U24").Value = 3
$ws.Range("V24").Value = 1.8
$ws.Range("W24").Value = 2
$ws.Range("X24").Value = 0.3999999999999999
$ws.Range("F27").Value = "NK Oriolik Oriovac"
$ws.Range("F28").Value = "NK Bistra"
$ws.Range("AA29").Value = -0.5
$ws.Range("AB29").Value = 0.45
$ws.Range("AC29").Value = 0.8999999999999999
$ws.Range("AD29").Value = -1
$ws.Range("B29").Value = 7250138
$ws.Range("E29").Value = "NK Tomislav"
$ws.Range("F29").Value = "Sloga Nova Gradiska"
$ws.Range("G29").Value = 2
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = "D"
$ws.Range("L29").Value = 2.1
$ws.Range("M29").Value = 3.4
$ws.Range("N29").Value = 2.9
$ws.Range("O29").Value = 2.1
$ws.Range("P29").Value = 3.4
$ws.Range("Q29").Value = 2.9
$ws.Range("R29").Value = -0.25
$ws.Range("S29").Value = 1.9
$ws.Range("T29").Value = 1.9
$ws.Range("U29").Value = 3
$ws.Range("V29").Value = 1.9
$ws.Range("W29").Value = 1.9
$ws.Range("Y29").Value = 2.4
$ws.Range("Z29").Value = -1
$ws.Range("AA30").Value = -1
$ws.Range("AB30").Value = 0.825
$ws.Range("AC30").Value = -1
$ws.Range("AD30").Value = 1
$ws.Range("B30").Value = 7250137
$ws.Range("E30").Value = "NK Granicar Zupanja"
$ws.Range("F30").Value = "NK Svacic"
$ws.Range("G30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1
$ws.Range("K30").Value = "A"
$ws.Range("L30").Value = 1.727
$ws.Range("M30").Value = 3.75
$ws.Range("N30").Value = 3.75
$ws.Range("O30").Value = 1.727
$ws.Range("P30").Value = 3.75
$ws.Range("Q30").Value = 3.75
$ws.Range("R30").Value = -0.75
$ws.Range("S30").Value = 1.975
$ws.Range("T30").Value = 1.825
$ws.Range("U30").Value = 2.5
$ws.Range("V30").Value = 1.8
$ws.Range("W30").Value = 2
$ws.Range("Y30").Value = -1
$ws.Range("Z30").Value = 2.75
$ws.Range("F34").Value = "Sava Strmec"
$ws.Range("E38").Value = "Sava Strmec"
$ws.Range("F38").Value = "NK Bistra"
$ws.Range("E43").Value = "NK Tomislav"
$ws.Range("F45").Value = "Sloga Nova Gradiska"
$ws.Range("E53").Value = "NK Tomislav"
$ws.Range("E54").Value = "NK Bistra"
$ws.Range("F68").Value = "Sloga Nova Gradiska"
$ws.Range("F70").Value = "NK Bistra"
$ws.Range("F71").Value = "NK Tomislav"
$ws.Range("E78").Value = "NK Bistra"
$ws.Range("E79").Value = "Sava Strmec"
$ws.Range("F84").Value = "Sava Strmec"
$ws.Range("E92").Value = "Sava Strmec"
$ws.Range("F94").Value = "NK Tomislav"
$ws.Range("F99").Value = "NK Vrapce"
$ws.Range("F100").Value = "NK Oriolik Oriovac"
$ws.Range("E101").Value = "Sava Strmec"
$ws.Range("E102").Value = "NK Bistra"
$ws.Range("AA107").Value = 0.825
$ws.Range("AD107").Value = 0.8500000000000001
$ws.Range("B107").Value = 8061518
$ws.Range("E107").Value = "NK Granicar Zupanja"
$ws.Range("F107").Value = "NK Croatia Dakovo"
$ws.Range("G107").Value = 1
$ws.Range("I107").Value = 1
$ws.Range("L107").Value = 2.5
$ws.Range("M107").Value = 3.4
$ws.Range("N107").Value = 2.375
$ws.Range("O107").Value = 2.375
$ws.Range("P107").Value = 3.4
$ws.Range("Q107").Value = 2.5
$ws.Range("R107").Value = 0
$ws.Range("S107").Value = 1.825
$ws.Range("T107").Value = 1.975
$ws.Range("U107").Value = 2.75
$ws.Range("V107").Value = 1.95
$ws.Range("W107").Value = 1.85
$ws.Range("X107").Value = 1.375
$ws.Range("AA108").Value = 0.95
$ws.Range("AD108").Value = 1
$ws.Range("B108").Value = 8061387
$ws.Range("E108").Value = "NK Neretvanac Opuzen"
$ws.Range("F108").Value = "RNK Split"
$ws.Range("G108").Value = 2
$ws.Range("I108").Value = 0
$ws.Range("L108").Value = 1.533
$ws.Range("M108").Value = 4
$ws.Range("N108").Value = 4.75
$ws.Range("O108").Value = 1.4
$ws.Range("P108").Value = 4.333
$ws.Range("Q108").Value = 5.75
$ws.Range("R108").Value = -1.25
$ws.Range("S108").Value = 1.95
$ws.Range("T108").Value = 1.85
$ws.Range("U108").Value = 3
$ws.Range("V108").Value = 1.8
$ws.Range("W108").Value = 2
$ws.Range("X108").Value = 0.3999999999999999
$ws.Range("E109").Value = "Sava Strmec"
$ws.Range("E110").Value = "NK Bistra"
$ws.Range("F119").Value = "Sloga Nova Gradiska"
$ws.Range("E120").Value = "NK Bistra"
$ws.Range("F120").Value = "Sava Strmec"
$ws.Range("E125").Value = "Sava Strmec"
$ws.Range("F127").Value = "Sava Strmec"
$ws.Range("E128").Value = "NK Bistra"
$ws.Range("E134").Value = "NK Bistra"
$ws.Range("E141").Value = "Sava Strmec"
$ws.Range("E142").Value = "NK Bistra"
$ws.Range("F145").Value = "NK Vrapce"
